$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Header row (row 1) ----
$ws.Range("A1").Value = "S.No."
$ws.Range("B1").Value = "Account"
$ws.Range("C1").Value = "Liability & Equity"
$ws.Range("D1").Value = "Assets"
$ws.Range("E1").Value = "Total"

# ---- Data rows (2-11): S.No., Account, Liability & Equity, Assets ----
$data = @(
    @(1, "Capital", 2000000, 0),
    @(2, "IOUiPersonX", 100000, 0),
    @(3, "InvPPFSBI", 0, 800000),
    @(4, "InvNPS", 0, 300000),
    @(5, "EquitySharesLarsenAndToubro", 0, 172181.05),
    @(6, "EquitySharesMandM", 0, 73100.01),
    @(7, "InvCoMfICICIPruDynamic", 0, 190000),
    @(8, "BaSBIAcc12345678", 0, 129680.62),
    @(9, "BaCitiAcc98765432", 0, 235038.32),
    @(10, "IOUiPersonY", 0, 200000)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $row++
}

# ---- Totals row (row 13; row 12 left blank, matching the template) ----
$ws.Range("B13").Value = "Total"
$ws.Range("C13").Value = 2100000
$ws.Range("D13").Value = 2100000
$ws.Range("E13").Value = 0

# ---- Column widths ----
$ws.Columns.Item(1).ColumnWidth = 6
$ws.Columns.Item(2).ColumnWidth = 45
$ws.Columns.Item(3).ColumnWidth = 12
$ws.Columns.Item(4).ColumnWidth = 12
$ws.Columns.Item(5).ColumnWidth = 12

# ---- Number formatting for monetary cells ----
$ws.Range("C2:D11").NumberFormat = "#,##0.00"
$ws.Range("C13:E13").NumberFormat = "#,##0.00"

# ---- Header style: bold white Calibri 11 text on a dark-blue fill ----
$headerRange = $ws.Range("A1:E1")
$headerRange.Font.Name = "Calibri"
$headerRange.Font.Size = 11
$headerRange.Font.Bold = $true
$headerRange.Font.Color = 16777215
$headerRange.Interior.Color = 9109504

# Row 1 is a bit taller to fit the larger header font
$ws.Rows.Item(1).RowHeight = 14.4

# ---- Selection, matching the saved view state ----
$ws.Range("E20").Select()
